# Generate Report for Archive
#
# Two files (4bf4c862-237b-4267-b283-1ead1eeac5c4.md and
# 6da2e667-b73e-414e-86e9-f846aa29f42b.md) moved from "Ready for handoff"
# to "In Translation" status. Update every sheet that tracks that status:
#   - Overview: columns "zh-cn" (E) and "de-de" (F), rows 3 and 4
#   - zh-cn:    column "Status" (C), rows 3 and 4
#   - de-de:    column "Status" (C), rows 3 and 4

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
